$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, pushing existing rows 180-219 down to 181-220
$ws.Rows.Item(180).Insert()

# After the insert, row 181 holds a copy of the data that was originally in row 180
# (rows below shifted down by one, retaining their formatting/values). Populate the
# newly blank row 180 with that same data, then overwrite the Fecha (D) and
# Volumen (J) columns with their new values.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(180, $col).Value2 = $ws.Cells.Item(181, $col).Value2
}

$ws.Cells.Item(180, 4).Value2 = 45275
$ws.Cells.Item(180, 10).Value2 = 90
